$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E30").NumberFormat = "@"

$ws.Range("E2").Value = "August/17/2003"
$ws.Range("E3").Value = "November/22/1989"
$ws.Range("E4").Value = "June/27/1943"
$ws.Range("E5").Value = "December/07/1991"
$ws.Range("E6").Value = "January/22/1961"
$ws.Range("E7").Value = "August/31/1971"
$ws.Range("E8").Value = "February/02/1957"
$ws.Range("E9").Value = "November/22/1957"
$ws.Range("E10").Value = "December/29/1970"
$ws.Range("E11").Value = "January/03/1988"
$ws.Range("E12").Value = "January/20/1986"
$ws.Range("E13").Value = "June/04/1971"
$ws.Range("E14").Value = "March/17/1994"
$ws.Range("E15").Value = "April/28/1988"
$ws.Range("E16").Value = "August/31/1979"
$ws.Range("E17").Value = "December/08/1972"
$ws.Range("E18").Value = "December/10/1935"
$ws.Range("E19").Value = "November/03/1948"
$ws.Range("E20").Value = "June/27/1926"
$ws.Range("E21").Value = "Jully/26/1960"
$ws.Range("E22").Value = "Jully/08/1937"
$ws.Range("E23").Value = "January/25/1963"
$ws.Range("E24").Value = "September/29/1974"
$ws.Range("E25").Value = "January/26/1962"
$ws.Range("E26").Value = "Jully/10/1950"
$ws.Range("E27").Value = "April/23/1969"
$ws.Range("E28").Value = "December/12/1978"
$ws.Range("E29").Value = "June/25/1972"
$ws.Range("E30").Value = "September/16/1998"

$ws.Columns("E").ColumnWidth = 22
$ws.Range("E10").Select()
